# "Added my part to the progress report"
#
# The author placed their cursor in the middle of the word "taking"
# (right after "us t") and kept typing / editing there, which is why
# Word's "_GoBack" bookmark (last-edit-location marker) ends up sitting
# between "...us t" and "aking...". That edit point splits the original
# run in two (the run boundary itself carries no visible text change -
# the sentence reads identically before and after), and the bookmark
# that used to sit at the end of the paragraph moves to the split point.

$d = $word.ActiveDocument

# Locate the sentence that contains the edit point.
$find = $d.Content.Find
$found = $find.Execute("us taking a major shift into the way", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target sentence in document"
}

# The split point sits right after "...us t" and before "aking...".
$sentenceStart = $find.Parent.Start
$splitPoint = $sentenceStart + 4

# Move the existing "_GoBack" bookmark (it marks the author's last edit
# location) from wherever it currently is to the split point.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$splitRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $splitRange)
